{"js": "const body = context.document.body;\n\n// 1) Touch the \"competitive for Redwood's agents...\" sentence span so the\n//    run split across several identically-formatted runs collapses into a\n//    single run (no visible text change \u2014 same characters before/after).\nconst introTarget =\n  \"competitive for Redwood\\u2019s agents, allowing them to better inform \" +\n  \"customers about trends based on information about the homes\";\nconst introResults = body.search(introTarget, { matchCase: true });\nintroResults.load(\"text\");\nawait context.sync();\n\nif (introResults.items.length > 0) {\n  introResults.items[0].insertText(introTarget, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the trailing \", Listing Date\" from the Facts line in the table\n//    (the actual semantic edit: \"Removed date from IP fact line\").\nconst factsResults = body.search(\", Listing Date\", { matchCase: true });\nfactsResults.load(\"text\");\nawait context.sync();\n\nif (factsResults.items.length > 0) {\n  factsResults.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Touch the \"competitive for Redwood's agents...\" sentence span so the\n#    run split across several identically-formatted runs collapses into a\n#    single run (no visible text change \u2014 same characters before/after).\n$introTarget = \"competitive for Redwood\" + [char]0x2019 + \"s agents, allowing them to better inform customers about trends based on information about the homes\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $introTarget\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\nif ($find.Execute()) {\n    $rng = $find.Parent\n    $rng.Text = $introTarget\n}\n\n# 2) Remove the trailing \", Listing Date\" from the Facts line in the table\n#    (the actual semantic edit: \"Removed date from IP fact line\").\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \", Listing Date\"\n$find2.MatchCase = $true\n$find2.Forward = $true\n$find2.Wrap = 0\n$find2.Replacement.Text = \"\"\n$find2.Execute($null, $true, $null, $null, $null, $null, $true, 0, $null, \"\", 2)\n"}
